$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) - update column F (想去人数) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 720
$ws1.Range("F3").Value = 507
$ws1.Range("F4").Value = 547
$ws1.Range("F7").Value = 36
$ws1.Range("F11").Value = 4627
$ws1.Range("F12").Value = 4431
$ws1.Range("F13").Value = 12

# Sheet "全部类型" (sheet4) - update column F (想去人数) values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 720
$ws4.Range("F3").Value = 507
$ws4.Range("F4").Value = 547
$ws4.Range("F7").Value = 36
$ws4.Range("F11").Value = 4628
$ws4.Range("F12").Value = 4431
$ws4.Range("F13").Value = 12
